$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player statistics cells to their new values.
# Values are written via the leading-apostrophe trick so Excel stores
# them as text (matching the workbook's inlineStr cell convention)
# rather than auto-converting numeric-looking text into a Number,
# then the cell style is reset to Normal so no stray NumberFormat/style
# gets attached to the cell.
$updates = @{
    "E4" = "90"
    "F4" = "1"
    "G4" = "1"
    "J5" = "1"
    "E7" = "360"
    "F7" = "4"
    "G7" = "4"
    "E8" = "113"
    "F8" = "4"
    "H8" = "4"
    "J8" = "4"
    "J10" = "3"
    "J11" = "2"
    "E12" = "360"
    "F12" = "4"
    "G12" = "4"
    "E13" = "346"
    "F13" = "4"
    "G13" = "4"
    "E14" = "247"
    "F14" = "4"
    "G14" = "4"
    "I14" = "4"
    "E15" = "290"
    "F15" = "4"
    "G15" = "4"
    "I15" = "3"
    "E17" = "99"
    "F17" = "4"
    "H17" = "4"
    "J17" = "4"
    "E18" = "292"
    "F18" = "4"
    "G18" = "4"
    "I18" = "4"
    "E19" = "215"
    "F19" = "4"
    "G19" = "2"
    "E20" = "69"
    "F20" = "3"
    "H20" = "3"
    "J20" = "4"
    "L20" = "1"
    "E21" = "11"
    "F21" = "2"
    "H21" = "2"
    "J21" = "4"
    "J23" = "2"
    "E25" = "356"
    "F25" = "4"
    "G25" = "4"
    "I25" = "1"
    "E26" = "352"
    "F26" = "4"
    "G26" = "4"
    "L26" = "1"
    "E27" = "303"
    "F27" = "4"
    "G27" = "4"
    "I27" = "3"
    "E28" = "37"
    "F28" = "2"
    "H28" = "2"
    "J28" = "3"
    "L28" = "1"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}

